# Update TSRs (Term Source Refs) in ENA "Single CDS annotated sequence" template.

$wb = $excel.ActiveWorkbook

# --- 1. Bump template version on the metadata sheet ---
$meta = $wb.Worksheets.Item("isa_template")
$meta.Range("B4").Value = "1.0.2"

# --- 2. Update header row + data row on the annotation table sheet ---
$ws = $wb.Worksheets.Item("New Table")

# Header row (row 1) - rename columns
$ws.Range("A1").Value = "Input [Data]"
$ws.Range("B1").Value = "Characteristic [organism]"
$ws.Range("C1").Value = "Term Source REF (OBI:0100026)"
$ws.Range("D1").Value = "Term Accession Number (OBI:0100026)"
$ws.Range("AF1").Value = "Output [Data]"

# Data row (row 2) - update values
$ws.Range("E2").Value = "No"
$ws.Range("F2").Value = "NCIT"
$ws.Range("G2").Value = "https://bioregistry.io/NCIT:C49487"

$ws.Range("M2").Value = "https://bioregistry.io/GO:0003968"

$ws.Range("W2").Value = "No"
$ws.Range("X2").Value = "NCIT"
$ws.Range("Y2").Value = "https://bioregistry.io/NCIT:C49487"

$ws.Range("Z2").Value = "No"
$ws.Range("AA2").Value = "NCIT"
$ws.Range("AB2").Value = "https://bioregistry.io/NCIT:C49487"
